$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Ligand-/Receptor-expressing cell counts (col E, K) from 1 to 3,
# and corresponding derived expression/specificity metrics (cols G-J, M-T),
# per Dr Hou advice.

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 98.91277700000001
$ws.Range("H2").Value = 296.738331
$ws.Range("I2").Value = 0.8120825131376513
$ws.Range("J2").Value = 0.8120825131376513
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 14.85604233333333
$ws.Range("N2").Value = 44.568127
$ws.Range("O2").Value = 0.09286934904108346
$ws.Range("P2").Value = 0.09286934904108346
$ws.Range("Q2").Value = 1469.45240241956
$ws.Range("R2").Value = 13225.07162177604
$ws.Range("S2").Value = 0.07541757436274078
$ws.Range("T2").Value = 0.07541757436274078

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 98.91277700000001
$ws.Range("H3").Value = 296.738331
$ws.Range("I3").Value = 0.8120825131376513
$ws.Range("J3").Value = 0.8120825131376513
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 114.2734143333333
$ws.Range("N3").Value = 342.820243
$ws.Range("O3").Value = 0.7143556381787382
$ws.Range("P3").Value = 0.7143556381787382
$ws.Range("Q3").Value = 11303.1007489816
$ws.Range("R3").Value = 101727.9067408344
$ws.Range("S3").Value = 0.5801157219262404
$ws.Range("T3").Value = 0.5801157219262404

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 98.91277700000001
$ws.Range("H4").Value = 296.738331
$ws.Range("I4").Value = 0.8120825131376513
$ws.Range("J4").Value = 0.8120825131376513
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 30.83766366666667
$ws.Range("N4").Value = 92.512991
$ws.Range("O4").Value = 0.1927750127801784
$ws.Range("P4").Value = 0.1927750127801784
$ws.Range("Q4").Value = 3050.238949462002
$ws.Range("R4").Value = 27452.15054515802
$ws.Range("S4").Value = 0.1565492168486701
$ws.Range("T4").Value = 0.1565492168486701

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 17.04862266666667
$ws.Range("H5").Value = 51.14586800000001
$ws.Range("I5").Value = 0.1399706767982279
$ws.Range("J5").Value = 0.1399706767982279
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 14.85604233333333
$ws.Range("N5").Value = 44.568127
$ws.Range("O5").Value = 0.09286934904108346
$ws.Range("P5").Value = 0.09286934904108346
$ws.Range("Q5").Value = 253.2750600610263
$ws.Range("R5").Value = 2279.475540549236
$ws.Range("S5").Value = 0.01299898563909131
$ws.Range("T5").Value = 0.01299898563909131

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 17.04862266666667
$ws.Range("H6").Value = 51.14586800000001
$ws.Range("I6").Value = 0.1399706767982279
$ws.Range("J6").Value = 0.1399706767982279
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 114.2734143333333
$ws.Range("N6").Value = 342.820243
$ws.Range("O6").Value = 0.7143556381787382
$ws.Range("P6").Value = 0.7143556381787382
$ws.Range("Q6").Value = 1948.204321800659
$ws.Range("R6").Value = 17533.83889620593
$ws.Range("S6").Value = 0.09998884215050802
$ws.Range("T6").Value = 0.09998884215050802

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 17.04862266666667
$ws.Range("H7").Value = 51.14586800000001
$ws.Range("I7").Value = 0.1399706767982279
$ws.Range("J7").Value = 0.1399706767982279
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 30.83766366666667
$ws.Range("N7").Value = 92.512991
$ws.Range("O7").Value = 0.1927750127801784
$ws.Range("P7").Value = 0.1927750127801784
$ws.Range("Q7").Value = 525.7396917745766
$ws.Range("R7").Value = 4731.657225971188
$ws.Range("S7").Value = 0.02698284900862861
$ws.Range("T7").Value = 0.02698284900862861

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.839988000000001
$ws.Range("H8").Value = 17.519964
$ws.Range("I8").Value = 0.0479468100641207
$ws.Range("J8").Value = 0.04794681006412069
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 14.85604233333333
$ws.Range("N8").Value = 44.568127
$ws.Range("O8").Value = 0.09286934904108346
$ws.Range("P8").Value = 0.09286934904108346
$ws.Range("Q8").Value = 86.75910895415868
$ws.Range("R8").Value = 780.831980587428
$ws.Range("S8").Value = 0.004452789039251359
$ws.Range("T8").Value = 0.004452789039251358

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.839988000000001
$ws.Range("H9").Value = 17.519964
$ws.Range("I9").Value = 0.0479468100641207
$ws.Range("J9").Value = 0.04794681006412069
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 114.2734143333333
$ws.Range("N9").Value = 342.820243
$ws.Range("O9").Value = 0.7143556381787382
$ws.Range("P9").Value = 0.7143556381787382
$ws.Range("Q9").Value = 667.3553684256948
$ws.Range("R9").Value = 6006.198315831252
$ws.Range("S9").Value = 0.03425107410198969
$ws.Range("T9").Value = 0.03425107410198969

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 5.839988000000001
$ws.Range("H10").Value = 17.519964
$ws.Range("I10").Value = 0.0479468100641207
$ws.Range("J10").Value = 0.04794681006412069
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 30.83766366666667
$ws.Range("N10").Value = 92.512991
$ws.Range("O10").Value = 0.1927750127801784
$ws.Range("P10").Value = 0.1927750127801784
$ws.Range("Q10").Value = 180.0915857613694
$ws.Range("R10").Value = 1620.824271852324
$ws.Range("S10").Value = 0.009242946922879655
$ws.Range("T10").Value = 0.009242946922879653

